# Reorder slides: swap the 4th and 5th slides in the deck.
# Before: ... , sldId=259 ("OOP Diagram") at position 4, sldId=265 ("App Tham Khao") at position 5
# After:  ... , sldId=265 ("App Tham Khao") at position 4, sldId=259 ("OOP Diagram") at position 5
$p = $ppt.ActivePresentation

# Move the slide currently in position 5 up to position 4 (pushes the former
# position-4 slide down to position 5), i.e. swap slides 4 and 5.
$p.Slides.Item(5).MoveTo(4)
